$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.694.35'
$ws.Range("E2").Value = '  +0.97%  '

$ws.Range("D3").Value = '1.645.89'
$ws.Range("E3").Value = '  +0.17%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.40%  '

$ws.Range("E6").Value = '  -1.03%  '

$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.21'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.65%  '

$ws.Range("E9").Value = '  +0.81%  '

$ws.Range("E10").Value = '  +0.51%  '

$ws.Range("E11").Value = '  +0.23%  '

$ws.Range("D12").Value = '1.879.84'
$ws.Range("E12").Value = '  +0.22%  '

$ws.Range("D13").Value = '1.654.19'
$ws.Range("E13").Value = '  +0.59%  '

$ws.Range("E14").Value = '  +0.25%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.562'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.84'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.87%  '

$ws.Range("D17").Value = '27.680.52'
$ws.Range("E17").Value = '  +1.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '231.56'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.34%  '

$ws.Range("D19").Value = '0.0₃0725'
$ws.Range("E19").Value = '  +0.77%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.64'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.20%  '

$ws.Range("E21").Value = '  +0.08%  '

$ws.Range("E22").Value = '  -0.57%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.86%  '

$ws.Range("E24").Value = '  -2.86%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.92'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.92'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.14%  '

$ws.Range("E27").Value = '  -1.64%  '

$ws.Range("B28").Value = 'BinanceUSD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.12%  '

$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.63'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.49%  '

$ws.Range("E30").Value = '  +0.52%  '

$ws.Range("E31").Value = '  +0.46%  '

$ws.Range("E32").Value = '  +0.66%  '

$ws.Range("D33").Value = '1.441.74'
$ws.Range("E33").Value = '  +2.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.14'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.14%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.59'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.44%  '

$ws.Range("E36").Value = '  -1.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.570'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.36%  '

$ws.Range("E38").Value = '  -0.15%  '

$ws.Range("E39").Value = '  +0.28%  '

$ws.Range("E40").Value = '  +11.94%  '

$ws.Range("E41").Value = '  -0.03%  '

$ws.Range("E42").Value = '  +0.07%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.61'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.62%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '67.51'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.37%  '

$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.26'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.94%  '

$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.789.09'
$ws.Range("E46").Value = '  +0.02%  '

$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.72'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.83%  '

$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.71'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.91%  '

$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0107'
$ws.Range("E49").Value = '  +2.21%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0987'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.04%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.97%  '
